$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation is inserted at row 90, pushing the existing
# rows (90..142) down by one (to 91..143).
$ws.Rows("90:90").Insert()

$ws.Cells.Item(90, 1).Value = 4
$ws.Cells.Item(90, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(90, 3).Value = "Los Lagos"
$ws.Cells.Item(90, 4).Value = 44460
$ws.Cells.Item(90, 5).Value = 10
$ws.Cells.Item(90, 6).Value = 100112021
$ws.Cells.Item(90, 7).Value = "Ají"
$ws.Cells.Item(90, 8).Value = "Inferno"
$ws.Cells.Item(90, 9).Value = "Segunda"
$ws.Cells.Item(90, 10).Value = 70
$ws.Cells.Item(90, 11).Value = 38000
$ws.Cells.Item(90, 12).Value = 38000
$ws.Cells.Item(90, 13).Value = 38000
$ws.Cells.Item(90, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(90, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(90, 16).Value = 3167
$ws.Cells.Item(90, 17).Value = 12
$ws.Cells.Item(90, 18).Value = "Hortaliza"
